$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 26,20
$arr[0,0] = "Sending cluster"
$arr[0,1] = "Ligand symbol"
$arr[0,2] = "Receptor symbol"
$arr[0,3] = "Target cluster"
$arr[0,4] = "Ligand-expressing cells"
$arr[0,5] = "Ligand detection rate"
$arr[0,6] = "Ligand average expression value"
$arr[0,7] = "Ligand total expression value"
$arr[0,8] = "Ligand derived specificity of average expression value"
$arr[0,9] = "Ligand derived specificity of total expression value"
$arr[0,10] = "Receptor-expressing cells"
$arr[0,11] = "Receptor detection rate"
$arr[0,12] = "Receptor average expression value"
$arr[0,13] = "Receptor total expression value"
$arr[0,14] = "Receptor derived specificity of average expression value"
$arr[0,15] = "Receptor derived specificity of total expression value"
$arr[0,16] = "Edge average expression weight"
$arr[0,17] = "Edge total expression weight"
$arr[0,18] = "Edge average expression derived specificity"
$arr[0,19] = "Edge total expression derived specificity"
$arr[1,0] = "ECs"
$arr[1,1] = "Fgf2"
$arr[1,2] = "Cd44"
$arr[1,3] = "ECs"
$arr[1,4] = 2
$arr[1,5] = 0.6666666666666666
$arr[1,6] = 1.066124666666667
$arr[1,7] = 3.198374
$arr[1,8] = 0.1044113535211941
$arr[1,9] = 0.1044113535211941
$arr[1,10] = 3
$arr[1,11] = 1
$arr[1,12] = 19.21315233333334
$arr[1,13] = 57.63945700000001
$arr[1,14] = 0.04451179209991234
$arr[1,15] = 0.04451179209991233
$arr[1,16] = 20.4836156269909
$arr[1,17] = 184.352540642918
$arr[1,18] = 0.004647536460805842
$arr[1,19] = 0.004647536460805841
$arr[2,0] = "ECs"
$arr[2,1] = "Fgf2"
$arr[2,2] = "Cd44"
$arr[2,3] = "FAPs"
$arr[2,4] = 2
$arr[2,5] = 0.6666666666666666
$arr[2,6] = 1.066124666666667
$arr[2,7] = 3.198374
$arr[2,8] = 0.1044113535211941
$arr[2,9] = 0.1044113535211941
$arr[2,10] = 3
$arr[2,11] = 1
$arr[2,12] = 92.44713066666667
$arr[2,13] = 277.341392
$arr[2,14] = 0.2141755495962477
$arr[2,15] = 0.2141755495962477
$arr[2,16] = 98.56016636628979
$arr[2,17] = 887.041497296608
$arr[2,18] = 0.02236235902448987
$arr[2,19] = 0.02236235902448986
$arr[3,0] = "ECs"
$arr[3,1] = "Fgf2"
$arr[3,2] = "Cd44"
$arr[3,3] = "Inflammatory-Mac"
$arr[3,4] = 2
$arr[3,5] = 0.6666666666666666
$arr[3,6] = 1.066124666666667
$arr[3,7] = 3.198374
$arr[3,8] = 0.1044113535211941
$arr[3,9] = 0.1044113535211941
$arr[3,10] = 3
$arr[3,11] = 1
$arr[3,12] = 166.8580016666666
$arr[3,13] = 500.5740049999999
$arr[3,14] = 0.3865658561145097
$arr[3,15] = 0.3865658561145097
$arr[3,16] = 177.8914314075411
$arr[3,17] = 1601.02288266787
$arr[3,18] = 0.04036186426199513
$arr[3,19] = 0.04036186426199513
$arr[4,0] = "ECs"
$arr[4,1] = "Fgf2"
$arr[4,2] = "Cd44"
$arr[4,3] = "MuSCs"
$arr[4,4] = 2
$arr[4,5] = 0.6666666666666666
$arr[4,6] = 1.066124666666667
$arr[4,7] = 3.198374
$arr[4,8] = 0.1044113535211941
$arr[4,9] = 0.1044113535211941
$arr[4,10] = 3
$arr[4,11] = 1
$arr[4,12] = 41.09915599999999
$arr[4,13] = 123.297468
$arr[4,14] = 0.09521587377309249
$arr[4,15] = 0.09521587377309249
$arr[4,16] = 43.81682399078133
$arr[4,17] = 394.351415917032
$arr[4,18] = 0.009941618257351754
$arr[4,19] = 0.009941618257351754
$arr[5,0] = "ECs"
$arr[5,1] = "Fgf2"
$arr[5,2] = "Cd44"
$arr[5,3] = "Resolving-Mac"
$arr[5,4] = 2
$arr[5,5] = 0.6666666666666666
$arr[5,6] = 1.066124666666667
$arr[5,7] = 3.198374
$arr[5,8] = 0.1044113535211941
$arr[5,9] = 0.1044113535211941
$arr[5,10] = 3
$arr[5,11] = 1
$arr[5,12] = 112.0244103333333
$arr[5,13] = 336.073231
$arr[5,14] = 0.2595309284162377
$arr[5,15] = 0.2595309284162377
$arr[5,16] = 119.4319871251549
$arr[5,17] = 1074.887884126394
$arr[5,18] = 0.02709797551655152
$arr[5,19] = 0.02709797551655152
$arr[6,0] = "FAPs"
$arr[6,1] = "Fgf2"
$arr[6,2] = "Cd44"
$arr[6,3] = "ECs"
$arr[6,4] = 3
$arr[6,5] = 1
$arr[6,6] = 7.418580000000001
$arr[6,7] = 22.25574
$arr[6,8] = 0.7265416542955204
$arr[6,9] = 0.7265416542955204
$arr[6,10] = 3
$arr[6,11] = 1
$arr[6,12] = 19.21315233333334
$arr[6,13] = 57.63945700000001
$arr[6,14] = 0.04451179209991234
$arr[6,15] = 0.04451179209991233
$arr[6,16] = 142.53430763702
$arr[6,17] = 1282.80876873318
$arr[6,18] = 0.03233967106792859
$arr[6,19] = 0.03233967106792858
$arr[7,0] = "FAPs"
$arr[7,1] = "Fgf2"
$arr[7,2] = "Cd44"
$arr[7,3] = "FAPs"
$arr[7,4] = 3
$arr[7,5] = 1
$arr[7,6] = 7.418580000000001
$arr[7,7] = 22.25574
$arr[7,8] = 0.7265416542955204
$arr[7,9] = 0.7265416542955204
$arr[7,10] = 3
$arr[7,11] = 1
$arr[7,12] = 92.44713066666667
$arr[7,13] = 277.341392
$arr[7,14] = 0.2141755495962477
$arr[7,15] = 0.2141755495962477
$arr[7,16] = 685.8264346211201
$arr[7,17] = 6172.43791159008
$arr[7,18] = 0.1556074581133101
$arr[7,19] = 0.1556074581133101
$arr[8,0] = "FAPs"
$arr[8,1] = "Fgf2"
$arr[8,2] = "Cd44"
$arr[8,3] = "Inflammatory-Mac"
$arr[8,4] = 3
$arr[8,5] = 1
$arr[8,6] = 7.418580000000001
$arr[8,7] = 22.25574
$arr[8,8] = 0.7265416542955204
$arr[8,9] = 0.7265416542955204
$arr[8,10] = 3
$arr[8,11] = 1
$arr[8,12] = 166.8580016666666
$arr[8,13] = 500.5740049999999
$arr[8,14] = 0.3865658561145097
$arr[8,15] = 0.3865658561145097
$arr[8,16] = 1237.8494340043
$arr[8,17] = 11140.6449060387
$arr[8,18] = 0.2808561965956
$arr[8,19] = 0.2808561965956
$arr[9,0] = "FAPs"
$arr[9,1] = "Fgf2"
$arr[9,2] = "Cd44"
$arr[9,3] = "MuSCs"
$arr[9,4] = 3
$arr[9,5] = 1
$arr[9,6] = 7.418580000000001
$arr[9,7] = 22.25574
$arr[9,8] = 0.7265416542955204
$arr[9,9] = 0.7265416542955204
$arr[9,10] = 3
$arr[9,11] = 1
$arr[9,12] = 41.09915599999999
$arr[9,13] = 123.297468
$arr[9,14] = 0.09521587377309249
$arr[9,15] = 0.09521587377309249
$arr[9,16] = 304.89737671848
$arr[9,17] = 2744.07639046632
$arr[9,18] = 0.06917829844629607
$arr[9,19] = 0.06917829844629607
$arr[10,0] = "FAPs"
$arr[10,1] = "Fgf2"
$arr[10,2] = "Cd44"
$arr[10,3] = "Resolving-Mac"
$arr[10,4] = 3
$arr[10,5] = 1
$arr[10,6] = 7.418580000000001
$arr[10,7] = 22.25574
$arr[10,8] = 0.7265416542955204
$arr[10,9] = 0.7265416542955204
$arr[10,10] = 3
$arr[10,11] = 1
$arr[10,12] = 112.0244103333333
$arr[10,13] = 336.073231
$arr[10,14] = 0.2595309284162377
$arr[10,15] = 0.2595309284162377
$arr[10,16] = 831.0620500106602
$arr[10,17] = 7479.558450095941
$arr[10,18] = 0.1885600300723856
$arr[10,19] = 0.1885600300723856
$arr[11,0] = "Inflammatory-Mac"
$arr[11,1] = "Fgf2"
$arr[11,2] = "Cd44"
$arr[11,3] = "ECs"
$arr[11,4] = 2
$arr[11,5] = 0.6666666666666666
$arr[11,6] = 0.4336433333333334
$arr[11,7] = 1.30093
$arr[11,8] = 0.04246903649677213
$arr[11,9] = 0.04246903649677213
$arr[11,10] = 3
$arr[11,11] = 1
$arr[11,12] = 19.21315233333334
$arr[11,13] = 57.63945700000001
$arr[11,14] = 0.04451179209991234
$arr[11,15] = 0.04451179209991233
$arr[11,16] = 8.331655421667779
$arr[11,17] = 74.98489879501001
$arr[11,18] = 0.001890372923227911
$arr[11,19] = 0.00189037292322791
$arr[12,0] = "Inflammatory-Mac"
$arr[12,1] = "Fgf2"
$arr[12,2] = "Cd44"
$arr[12,3] = "FAPs"
$arr[12,4] = 2
$arr[12,5] = 0.6666666666666666
$arr[12,6] = 0.4336433333333334
$arr[12,7] = 1.30093
$arr[12,8] = 0.04246903649677213
$arr[12,9] = 0.04246903649677213
$arr[12,10] = 3
$arr[12,11] = 1
$arr[12,12] = 92.44713066666667
$arr[12,13] = 277.341392
$arr[12,14] = 0.2141755495962477
$arr[12,15] = 0.2141755495962477
$arr[12,16] = 40.08908189939556
$arr[12,17] = 360.80173709456
$arr[12,18] = 0.009095829232519274
$arr[12,19] = 0.009095829232519274
$arr[13,0] = "Inflammatory-Mac"
$arr[13,1] = "Fgf2"
$arr[13,2] = "Cd44"
$arr[13,3] = "Inflammatory-Mac"
$arr[13,4] = 2
$arr[13,5] = 0.6666666666666666
$arr[13,6] = 0.4336433333333334
$arr[13,7] = 1.30093
$arr[13,8] = 0.04246903649677213
$arr[13,9] = 0.04246903649677213
$arr[13,10] = 3
$arr[13,11] = 1
$arr[13,12] = 166.8580016666666
$arr[13,13] = 500.5740049999999
$arr[13,14] = 0.3865658561145097
$arr[13,15] = 0.3865658561145097
$arr[13,16] = 72.35686003607222
$arr[13,17] = 651.21174032465
$arr[13,18] = 0.01641707945173308
$arr[13,19] = 0.01641707945173308
$arr[14,0] = "Inflammatory-Mac"
$arr[14,1] = "Fgf2"
$arr[14,2] = "Cd44"
$arr[14,3] = "MuSCs"
$arr[14,4] = 2
$arr[14,5] = 0.6666666666666666
$arr[14,6] = 0.4336433333333334
$arr[14,7] = 1.30093
$arr[14,8] = 0.04246903649677213
$arr[14,9] = 0.04246903649677213
$arr[14,10] = 3
$arr[14,11] = 1
$arr[14,12] = 41.09915599999999
$arr[14,13] = 123.297468
$arr[14,14] = 0.09521587377309249
$arr[14,15] = 0.09521587377309249
$arr[14,16] = 17.82237500502666
$arr[14,17] = 160.40137504524
$arr[14,18] = 0.004043726418341513
$arr[14,19] = 0.004043726418341513
$arr[15,0] = "Inflammatory-Mac"
$arr[15,1] = "Fgf2"
$arr[15,2] = "Cd44"
$arr[15,3] = "Resolving-Mac"
$arr[15,4] = 2
$arr[15,5] = 0.6666666666666666
$arr[15,6] = 0.4336433333333334
$arr[15,7] = 1.30093
$arr[15,8] = 0.04246903649677213
$arr[15,9] = 0.04246903649677213
$arr[15,10] = 3
$arr[15,11] = 1
$arr[15,12] = 112.0244103333333
$arr[15,13] = 336.073231
$arr[15,14] = 0.2595309284162377
$arr[15,15] = 0.2595309284162377
$arr[15,16] = 48.57863871164778
$arr[15,17] = 437.2077484048301
$arr[15,18] = 0.01102202847095035
$arr[15,19] = 0.01102202847095035
$arr[16,0] = "MuSCs"
$arr[16,1] = "Fgf2"
$arr[16,2] = "Cd44"
$arr[16,3] = "ECs"
$arr[16,4] = 3
$arr[16,5] = 1
$arr[16,6] = 0.951285
$arr[16,7] = 2.853855
$arr[16,8] = 0.09316448398568379
$arr[16,9] = 0.09316448398568379
$arr[16,10] = 3
$arr[16,11] = 1
$arr[16,12] = 19.21315233333334
$arr[16,13] = 57.63945700000001
$arr[16,14] = 0.04451179209991234
$arr[16,15] = 0.04451179209991233
$arr[16,16] = 18.27718361741501
$arr[16,17] = 164.494652556735
$arr[16,18] = 0.004146918142266369
$arr[16,19] = 0.004146918142266368
$arr[17,0] = "MuSCs"
$arr[17,1] = "Fgf2"
$arr[17,2] = "Cd44"
$arr[17,3] = "FAPs"
$arr[17,4] = 3
$arr[17,5] = 1
$arr[17,6] = 0.951285
$arr[17,7] = 2.853855
$arr[17,8] = 0.09316448398568379
$arr[17,9] = 0.09316448398568379
$arr[17,10] = 3
$arr[17,11] = 1
$arr[17,12] = 92.44713066666667
$arr[17,13] = 277.341392
$arr[17,14] = 0.2141755495962477
$arr[17,15] = 0.2141755495962477
$arr[17,16] = 87.94356869624001
$arr[17,17] = 791.49211826616
$arr[17,18] = 0.01995355456048465
$arr[17,19] = 0.01995355456048464
$arr[18,0] = "MuSCs"
$arr[18,1] = "Fgf2"
$arr[18,2] = "Cd44"
$arr[18,3] = "Inflammatory-Mac"
$arr[18,4] = 3
$arr[18,5] = 1
$arr[18,6] = 0.951285
$arr[18,7] = 2.853855
$arr[18,8] = 0.09316448398568379
$arr[18,9] = 0.09316448398568379
$arr[18,10] = 3
$arr[18,11] = 1
$arr[18,12] = 166.8580016666666
$arr[18,13] = 500.5740049999999
$arr[18,14] = 0.3865658561145097
$arr[18,15] = 0.3865658561145097
$arr[18,16] = 158.729514115475
$arr[18,17] = 1428.565627039275
$arr[18,18] = 0.03601420851139239
$arr[18,19] = 0.03601420851139239
$arr[19,0] = "MuSCs"
$arr[19,1] = "Fgf2"
$arr[19,2] = "Cd44"
$arr[19,3] = "MuSCs"
$arr[19,4] = 3
$arr[19,5] = 1
$arr[19,6] = 0.951285
$arr[19,7] = 2.853855
$arr[19,8] = 0.09316448398568379
$arr[19,9] = 0.09316448398568379
$arr[19,10] = 3
$arr[19,11] = 1
$arr[19,12] = 41.09915599999999
$arr[19,13] = 123.297468
$arr[19,14] = 0.09521587377309249
$arr[19,15] = 0.09521587377309249
$arr[19,16] = 39.09701061545999
$arr[19,17] = 351.87309553914
$arr[19,18] = 0.008870737747316164
$arr[19,19] = 0.008870737747316164
$arr[20,0] = "MuSCs"
$arr[20,1] = "Fgf2"
$arr[20,2] = "Cd44"
$arr[20,3] = "Resolving-Mac"
$arr[20,4] = 3
$arr[20,5] = 1
$arr[20,6] = 0.951285
$arr[20,7] = 2.853855
$arr[20,8] = 0.09316448398568379
$arr[20,9] = 0.09316448398568379
$arr[20,10] = 3
$arr[20,11] = 1
$arr[20,12] = 112.0244103333333
$arr[20,13] = 336.073231
$arr[20,14] = 0.2595309284162377
$arr[20,15] = 0.2595309284162377
$arr[20,16] = 106.567141183945
$arr[20,17] = 959.1042706555052
$arr[20,18] = 0.02417906502422422
$arr[20,19] = 0.02417906502422422
$arr[21,0] = "Resolving-Mac"
$arr[21,1] = "Fgf2"
$arr[21,2] = "Cd44"
$arr[21,3] = "ECs"
$arr[21,4] = 3
$arr[21,5] = 1
$arr[21,6] = 0.3411786666666667
$arr[21,7] = 1.023536
$arr[21,8] = 0.03341347170082953
$arr[21,9] = 0.03341347170082953
$arr[21,10] = 3
$arr[21,11] = 1
$arr[21,12] = 19.21315233333334
$arr[21,13] = 57.63945700000001
$arr[21,14] = 0.04451179209991234
$arr[21,15] = 0.04451179209991233
$arr[21,16] = 6.555117695550224
$arr[21,17] = 58.99605925995201
$arr[21,18] = 0.001487293505683628
$arr[21,19] = 0.001487293505683628
$arr[22,0] = "Resolving-Mac"
$arr[22,1] = "Fgf2"
$arr[22,2] = "Cd44"
$arr[22,3] = "FAPs"
$arr[22,4] = 3
$arr[22,5] = 1
$arr[22,6] = 0.3411786666666667
$arr[22,7] = 1.023536
$arr[22,8] = 0.03341347170082953
$arr[22,9] = 0.03341347170082953
$arr[22,10] = 3
$arr[22,11] = 1
$arr[22,12] = 92.44713066666667
$arr[22,13] = 277.341392
$arr[22,14] = 0.2141755495962477
$arr[22,15] = 0.2141755495962477
$arr[22,16] = 31.54098877801245
$arr[22,17] = 283.868899002112
$arr[22,18] = 0.007156348665443834
$arr[22,19] = 0.007156348665443834
$arr[23,0] = "Resolving-Mac"
$arr[23,1] = "Fgf2"
$arr[23,2] = "Cd44"
$arr[23,3] = "Inflammatory-Mac"
$arr[23,4] = 3
$arr[23,5] = 1
$arr[23,6] = 0.3411786666666667
$arr[23,7] = 1.023536
$arr[23,8] = 0.03341347170082953
$arr[23,9] = 0.03341347170082953
$arr[23,10] = 3
$arr[23,11] = 1
$arr[23,12] = 166.8580016666666
$arr[23,13] = 500.5740049999999
$arr[23,14] = 0.3865658561145097
$arr[23,15] = 0.3865658561145097
$arr[23,16] = 56.92839053129777
$arr[23,17] = 512.35551478168
$arr[23,18] = 0.01291650729378911
$arr[23,19] = 0.01291650729378911
$arr[24,0] = "Resolving-Mac"
$arr[24,1] = "Fgf2"
$arr[24,2] = "Cd44"
$arr[24,3] = "MuSCs"
$arr[24,4] = 3
$arr[24,5] = 1
$arr[24,6] = 0.3411786666666667
$arr[24,7] = 1.023536
$arr[24,8] = 0.03341347170082953
$arr[24,9] = 0.03341347170082953
$arr[24,10] = 3
$arr[24,11] = 1
$arr[24,12] = 41.09915599999999
$arr[24,13] = 123.297468
$arr[24,14] = 0.09521587377309249
$arr[24,15] = 0.09521587377309249
$arr[24,16] = 14.02215524520533
$arr[24,17] = 126.199397206848
$arr[24,18] = 0.003181492903786982
$arr[24,19] = 0.003181492903786982
$arr[25,0] = "Resolving-Mac"
$arr[25,1] = "Fgf2"
$arr[25,2] = "Cd44"
$arr[25,3] = "Resolving-Mac"
$arr[25,4] = 3
$arr[25,5] = 1
$arr[25,6] = 0.3411786666666667
$arr[25,7] = 1.023536
$arr[25,8] = 0.03341347170082953
$arr[25,9] = 0.03341347170082953
$arr[25,10] = 3
$arr[25,11] = 1
$arr[25,12] = 112.0244103333333
$arr[25,13] = 336.073231
$arr[25,14] = 0.2595309284162377
$arr[25,15] = 0.2595309284162377
$arr[25,16] = 38.22033895164623
$arr[25,17] = 343.983050564816
$arr[25,18] = 0.008671829332125972
$arr[25,19] = 0.008671829332125972
$ws.Range("A1:T26").Value = $arr
